$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 1839
$ws.Range("B4").Value = 1725
$ws.Range("B5").Value = 1497
$ws.Range("B6").Value = 1302
$ws.Range("B7").Value = 1018
$ws.Range("B8").Value = 850
$ws.Range("B9").Value = 615
$ws.Range("B10").Value = 495
$ws.Range("B11").Value = 451
$ws.Range("B12").Value = 321
